# Update cell data to reflect corrected "Fallimenti per Fragilita'" (F)
# and "Fallimenti per Obsolescenza" (E) counts for the "prova" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prova")

# --- "LLM" block (rows 4-9) ---
# Row 5
$ws.Range("E5").Value = 12
$ws.Range("F5").Value = 2

# Row 6
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 2

# Row 7
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 2

# Row 8
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 2

# Row 9
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 2

# --- "Analitica" block (rows 11-16) ---
# Row 12
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 2

# Row 13
$ws.Range("E13").Value = 5
$ws.Range("F13").Value = 2

# Row 14
$ws.Range("E14").Value = 6
$ws.Range("F14").Value = 2

# Row 15
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 2

# Row 16
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 2
